# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.748.61'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '1.889.41'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7903'
$ws.Range("E5").Value = '  -5.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.43'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3152'
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.44'
$ws.Range("E9").Value = '  -4.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06986'
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08038'
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7574'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").Value = '1.910.33'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.283'
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.15'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("D16").Value = '29.761.61'
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.77'
$ws.Range("E17").Value = '  -2.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.913'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.29'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007656'
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = '2.149.65'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.099'
$ws.Range("E23").Value = '  +16.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.274'
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.42'
$ws.Range("E27").Value = '  -3.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.60'
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.044'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.379'
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.533'
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.378'
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05674'
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.048'
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.258'
$ws.Range("E35").Value = '  -1.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7319'
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9985'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.573'
$ws.Range("E38").Value = '  -4.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01897'
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.771'
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4380'
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.19'
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.808'
$ws.Range("E43").Value = '  -2.86%  '
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8375'
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.19'
$ws.Range("E46").Value = '  +1.35%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.017.59'
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.846'
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.837'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.417'
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("D51").Value = '2.048.76'
$ws.Range("E51").Value = '  -0.65%  '
